{"js": "// Update the date heading paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2024-05-19 Sunday\", Word.InsertLocation.replace);\n\n// Update every cell in the addition/subtraction practice table with the\n// new set of problems (row-major order, same 20x5 grid as the original).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst newValues = [\n  [\"74-60=14\", \"66-29=37\", \"94-41=53\", \"75+16=91\", \"84-51=33\"],\n  [\"19+32=51\", \"4+87=91\", \"32+32=64\", \"30-5=25\", \"58+39=97\"],\n  [\"57+39=96\", \"2+23=25\", \"31-27=4\", \"17+28=45\", \"4+37=41\"],\n  [\"87-1=86\", \"53-8=45\", \"81+3=84\", \"37-0=37\", \"8+67=75\"],\n  [\"53+8=61\", \"63-61=2\", \"42-23=19\", \"57+26=83\", \"15+25=40\"],\n  [\"7+47=54\", \"26+48=74\", \"89-43=46\", \"34+40=74\", \"39+37=76\"],\n  [\"9+87=96\", \"77-8=69\", \"2+50=52\", \"16+44=60\", \"32+23=55\"],\n  [\"79-8=71\", \"79-18=61\", \"89-3=86\", \"74+13=87\", \"2+78=80\"],\n  [\"62+27=89\", \"13+52=65\", \"49-6=43\", \"46+20=66\", \"92+4=96\"],\n  [\"65+10=75\", \"83-8=75\", \"69-67=2\", \"8+10=18\", \"83-35=48\"],\n  [\"47+52=99\", \"90-13=77\", \"35+28=63\", \"88-63=25\", \"88-77=11\"],\n  [\"4+16=20\", \"44+34=78\", \"80-35=45\", \"49-5=44\", \"68-5=63\"],\n  [\"48+4=52\", \"77-7=70\", \"75-48=27\", \"6+26=32\", \"49-7=42\"],\n  [\"46+43=89\", \"36+50=86\", \"68+4=72\", \"66-47=19\", \"40-0=40\"],\n  [\"92-84=8\", \"42-19=23\", \"0+45=45\", \"20+60=80\", \"62-10=52\"],\n  [\"66+27=93\", \"51+5=56\", \"38-27=11\", \"24+34=58\", \"59-45=14\"],\n  [\"31-18=13\", \"58+40=98\", \"78-7=71\", \"35-30=5\", \"57+41=98\"],\n  [\"13+31=44\", \"39+24=63\", \"4+72=76\", \"82-78=4\", \"4+42=46\"],\n  [\"78+0=78\", \"95-43=52\", \"63+28=91\", \"9+16=25\", \"71+9=80\"],\n  [\"3+7=10\", \"37+38=75\", \"37+48=85\", \"72-28=44\", \"92-88=4\"]\n];\ntable.values = newValues;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph).\n$d.Paragraphs.Item(1).Range.Text = \"2024-05-19 Sunday\"\n\n# Update every cell in the addition/subtraction practice table.\n$t = $d.Tables.Item(1)\n$values = @(\n    @(\"74-60=14\", \"66-29=37\", \"94-41=53\", \"75+16=91\", \"84-51=33\"),\n    @(\"19+32=51\", \"4+87=91\", \"32+32=64\", \"30-5=25\", \"58+39=97\"),\n    @(\"57+39=96\", \"2+23=25\", \"31-27=4\", \"17+28=45\", \"4+37=41\"),\n    @(\"87-1=86\", \"53-8=45\", \"81+3=84\", \"37-0=37\", \"8+67=75\"),\n    @(\"53+8=61\", \"63-61=2\", \"42-23=19\", \"57+26=83\", \"15+25=40\"),\n    @(\"7+47=54\", \"26+48=74\", \"89-43=46\", \"34+40=74\", \"39+37=76\"),\n    @(\"9+87=96\", \"77-8=69\", \"2+50=52\", \"16+44=60\", \"32+23=55\"),\n    @(\"79-8=71\", \"79-18=61\", \"89-3=86\", \"74+13=87\", \"2+78=80\"),\n    @(\"62+27=89\", \"13+52=65\", \"49-6=43\", \"46+20=66\", \"92+4=96\"),\n    @(\"65+10=75\", \"83-8=75\", \"69-67=2\", \"8+10=18\", \"83-35=48\"),\n    @(\"47+52=99\", \"90-13=77\", \"35+28=63\", \"88-63=25\", \"88-77=11\"),\n    @(\"4+16=20\", \"44+34=78\", \"80-35=45\", \"49-5=44\", \"68-5=63\"),\n    @(\"48+4=52\", \"77-7=70\", \"75-48=27\", \"6+26=32\", \"49-7=42\"),\n    @(\"46+43=89\", \"36+50=86\", \"68+4=72\", \"66-47=19\", \"40-0=40\"),\n    @(\"92-84=8\", \"42-19=23\", \"0+45=45\", \"20+60=80\", \"62-10=52\"),\n    @(\"66+27=93\", \"51+5=56\", \"38-27=11\", \"24+34=58\", \"59-45=14\"),\n    @(\"31-18=13\", \"58+40=98\", \"78-7=71\", \"35-30=5\", \"57+41=98\"),\n    @(\"13+31=44\", \"39+24=63\", \"4+72=76\", \"82-78=4\", \"4+42=46\"),\n    @(\"78+0=78\", \"95-43=52\", \"63+28=91\", \"9+16=25\", \"71+9=80\"),\n    @(\"3+7=10\", \"37+38=75\", \"37+48=85\", \"72-28=44\", \"92-88=4\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $values[$r - 1][$c - 1]\n    }\n}\n\n"}
